$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextValue 'D2' '26.861.88'
$ws.Range('E2').Value = '  -2.14%  '
Set-TextValue 'D3' '1.565.51'
$ws.Range('E3').Value = '  -0.39%  '
$ws.Range('E4').Value = '  +0.13%  '
Set-TextValue 'D5' '206.33'
$ws.Range('E5').Value = '  -1.08%  '
Set-TextValue 'D6' '0.489'
$ws.Range('E6').Value = '  -2.03%  '
$ws.Range('E7').Value = '  +0.13%  '
Set-TextValue 'D8' '21.93'
$ws.Range('E8').Value = '  -1.62%  '
$ws.Range('E9').Value = '  -0.50%  '
$ws.Range('E10').Value = '  -1.72%  '
Set-TextValue 'D11' '0.0864'
$ws.Range('E11').Value = '  -0.37%  '
Set-TextValue 'D12' '1.789.40'
$ws.Range('E12').Value = '  -0.30%  '
Set-TextValue 'D13' '1.569.77'
$ws.Range('E13').Value = '  -0.17%  '
$ws.Range('E14').Value = '  -2.05%  '
$ws.Range('E15').Value = '  -1.19%  '
Set-TextValue 'D16' '26.874.24'
$ws.Range('E16').Value = '  -2.03%  '
Set-TextValue 'D17' '61.44'
$ws.Range('E17').Value = '  -3.43%  '
$ws.Range('B18').Value = 'BitcoinCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D18' '215.22'
$ws.Range('E18').Value = '  +1.00%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D19' '7.38'
$ws.Range('E19').Value = '  +0.80%  '
$ws.Range('E20').Value = '  -2.36%  '
$ws.Range('E21').Value = '  +0.16%  '
Set-TextValue 'D22' '4.14'
$ws.Range('E22').Value = '  +0.15%  '
Set-TextValue 'D23' '9.29'
$ws.Range('E23').Value = '  -2.81%  '
Set-TextValue 'D24' '2.00'
$ws.Range('E24').Value = '  -0.43%  '
Set-TextValue 'D25' '153.71'
$ws.Range('E25').Value = '  +0.22%  '
Set-TextValue 'D26' '6.70'
$ws.Range('E26').Value = '  +0.10%  '
Set-TextValue 'D27' '14.97'
$ws.Range('E27').Value = '  -0.19%  '
$ws.Range('E28').Value = '  +0.14%  '
Set-TextValue 'D29' '0.103'
$ws.Range('E29').Value = '  -1.31%  '
Set-TextValue 'D30' '0.0466'
$ws.Range('E30').Value = '  -1.33%  '
Set-TextValue 'D31' '1.11'
$ws.Range('E31').Value = '  -3.46%  '
Set-TextValue 'D32' '3.17'
$ws.Range('E32').Value = '  -1.00%  '
Set-TextValue 'D33' '1.402.71'
$ws.Range('E33').Value = '  +0.95%  '
$ws.Range('E34').Value = '  -1.83%  '
$ws.Range('E35').Value = '  -2.39%  '
$ws.Range('E36').Value = '  -0.65%  '
Set-TextValue 'D37' '0.933'
$ws.Range('E37').Value = '  -1.68%  '
Set-TextValue 'D38' '0.0162'
$ws.Range('E38').Value = '  -2.77%  '
Set-TextValue 'D39' '0.529'
$ws.Range('E39').Value = '  -1.29%  '
Set-TextValue 'D40' '0.813'
$ws.Range('E40').Value = '  -1.52%  '
$ws.Range('E41').Value = '  +0.12%  '
Set-TextValue 'D42' '0.988'
$ws.Range('E42').Value = '  -0.30%  '
Set-TextValue 'D43' '1.80'
$ws.Range('E43').Value = '  -0.04%  '
$ws.Range('E44').Value = '  +0.91%  '
Set-TextValue 'D45' '2.19'
$ws.Range('E45').Value = '  +0.58%  '
Set-TextValue 'D46' '63.14'
$ws.Range('E46').Value = '  -1.74%  '
Set-TextValue 'D47' '1.701.87'
$ws.Range('E47').Value = '  -0.19%  '
Set-TextValue 'D48' '86.32'
$ws.Range('E48').Value = '  +0.63%  '
Set-TextValue 'D49' '0.0₇0983'
$ws.Range('E49').Value = '  -0.84%  '
Set-TextValue 'D50' '0.0949'
$ws.Range('E50').Value = '  -0.67%  '
$ws.Range('E51').Value = '  -0.34%  '
